$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, shifting existing rows 15:111 down to 16:112
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with its values (mirrors the row that
# was previously at 15, but with updated Fecha/Volumen/Precio promedio
# ponderado/Precio $/Kg figures)
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 45168
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100112035
$ws.Range("G15").Value = "Bruselas (repollito)"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 360
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17361
$ws.Range("N15").Value = "`$/malla 15 kilos"
$ws.Range("O15").Value = "Provincia de Quillota"
$ws.Range("P15").Value = 1157
$ws.Range("Q15").Value = 15
$ws.Range("R15").Value = "Hortaliza"
